$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.006.54'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '1.560.96'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = "'208.26"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').Value = "'0.490"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').Value = "'22.08"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = "'0.249"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('D11').Value = "'0.0855"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '1.782.05'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '1.561.77'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').Value = "'3.73"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('E15').Value = '  +0.18%  '
$ws.Range('D16').Value = '27.021.20'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('D17').Value = "'61.83"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '0.0₃0706'
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').Value = "'215.53"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').Value = "'7.39"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('D23').Value = "'9.21"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('D24').Value = "'1.95"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').Value = "'153.14"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').Value = "'15.05"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('D28').Value = "'0.106"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.69%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('E31').Value = '  +3.56%  '
$ws.Range('D32').Value = "'3.24"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').Value = "'3.18"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.74%  '
$ws.Range('D34').Value = '1.423.13'
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('D35').Value = "'1.07"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +10.75%  '
$ws.Range('D36').Value = "'1.61"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.73%  '
$ws.Range('D37').Value = "'2.34"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.50%  '
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'5.80"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = "'0.810"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('D42').Value = "'1.01"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = "'2.31"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = "'0.999"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = "'64.69"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('E46').Value = '  -1.13%  '
$ws.Range('D47').Value = '1.696.46'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('D48').Value = "'86.77"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('E51').Value = '  +0.68%  '
